$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.343.07'
$ws.Range("E2").Value = '  +0.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.641.88'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '597.55'
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.56'
$ws.Range("E6").Value = '  +0.73%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.641.08'
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").Value = '  +7.88%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.353'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.18'
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000193'
$ws.Range("E15").Value = '  +2.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.123.32'
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '68.252.08'
$ws.Range("E17").Value = '  +0.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.645.81'
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.42'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '365.00'
$ws.Range("E20").Value = '  -2.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.50'
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  +3.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.90'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.57'
$ws.Range("E25").Value = '  +2.92%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.84'
$ws.Range("E27").Value = '  -1.38%  '
$ws.Range("E28").Value = '  +1.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.781.19'
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '573.13'
$ws.Range("E31").Value = '  -0.83%  '
$ws.Range("E32").Value = '  +4.14%  '
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("E34").Value = '  +0.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("E35").Value = '  +3.07%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.59'
$ws.Range("E37").Value = '  +5.29%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '161.02'
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.40'
$ws.Range("E39").Value = '  +1.34%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.90'
$ws.Range("E40").Value = '  +0.30%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.374'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.42'
$ws.Range("E42").Value = '  +1.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0338'
$ws.Range("E43").Value = '  +5.61%  '
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("E45").Value = '  +3.67%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.65'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("E47").Value = '  +0.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '156.75'
$ws.Range("E48").Value = '  +1.04%  '
$ws.Range("E49").Value = '  +1.84%  '
$ws.Range("E50").Value = '  +1.08%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.92'
$ws.Range("E51").Value = '  -0.05%  '
